$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("GeneralVariables")
$ws.Activate()

# Insert a new row before row 15 (pushes "idTestingCompanySOI66" block and
# everything after it down by one) and populate it with the new
# testingCompanySOI880 name/value pair.
$ws.Rows.Item(15).Insert()
$ws.Range("A15").Value = "testingCompanySOI880"
$ws.Range("B15").Value = "AutoTestingCompany_SOI880"

# Insert a second new row right before the "testingCompanyD02Functional"
# row (originally row 28, now row 29 after the first insertion above) and
# populate it with the matching id/value pair for the new company.
$ws.Rows.Item(29).Insert()
$ws.Range("A29").Value = "idTestingCompanySOI880"
$ws.Range("B29").Value = "0013E00001ABencQAD"

# Match the author's final selection in the workbook.
$ws.Range("F23").Select()
